# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect the newly generated gh-pages scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 394
$wsExhibit.Range("F3").Value = 120
$wsExhibit.Range("F5").Value = 0
$wsExhibit.Range("F7").Value = 0
$wsExhibit.Range("F8").Value = 147
$wsExhibit.Range("F9").Value = 65
$wsExhibit.Range("F10").Value = 517

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1641
$wsAll.Range("F6").Value = 23
$wsAll.Range("F7").Value = 419
$wsAll.Range("F9").Value = 0
$wsAll.Range("F10").Value = 517
